# Update countries & provincias Spain
# Refreshes the COVID-19 "Pais" sheet: new snapshot timestamp, updated
# case totals for several countries, and a handful of countries that
# leap-frogged their neighbour in the ranking (so the two rows swap
# which country name they display).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 3 de Julio de 2020 a las 09:33"

# --- Countries that swapped rank with their neighbour ------------------
# (identical/near-identical totals, so the row keeps its numbers but the
# country label flips with the adjacent row)
$ws.Range("A52").Value = "Armenia"
$ws.Range("A53").Value = "Nigeria"

$ws.Range("A203").Value = "Santa Lucia"
$ws.Range("A204").Value = "Laos"

$ws.Range("A205").Value = "Fiyi"
$ws.Range("A206").Value = "Dominica"

$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"

# --- Updated case numbers (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes) ------------------------------

# Row 30 - Belgica
$ws.Range("B30").Value = 61727
$ws.Range("C30").Value = 129
$ws.Range("D30").Value = 17073
$ws.Range("E30").Value = 34889
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = 9765

# Row 37 - Ucrania
$ws.Range("B37").Value = 46763
$ws.Range("C37").Value = 876
$ws.Range("D37").Value = 20558
$ws.Range("E37").Value = 24993
$ws.Range("G37").Value = 27
$ws.Range("H37").Value = 1212

# Row 38 - Singapur
$ws.Range("B38").Value = 44479
$ws.Range("C38").Value = 169
$ws.Range("E38").Value = 5024

# Row 52 - Armenia (now ranked above Nigeria)
$ws.Range("B52").Value = 27320
$ws.Range("C52").Value = 662
$ws.Range("D52").Value = 15484
$ws.Range("E52").Value = 11367
$ws.Range("G52").Value = 10
$ws.Range("H52").Value = 469

# Row 53 - Nigeria (now ranked below Armenia)
$ws.Range("B53").Value = 27110
$ws.Range("D53").Value = 10801
$ws.Range("E53").Value = 15693
$ws.Range("H53").Value = 616

# Row 71 - Sudan
$ws.Range("B71").Value = 9663
$ws.Range("C71").Value = 90
$ws.Range("D71").Value = 4624
$ws.Range("E71").Value = 4435
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 604

# Row 79 - El Salvador
$ws.Range("D79").Value = 4137
$ws.Range("E79").Value = 2663
$ws.Range("G79").Value = 9
$ws.Range("H79").Value = 200

# Row 95 - Hungria
$ws.Range("B95").Value = 4172
$ws.Range("C95").Value = 6
$ws.Range("D95").Value = 2752
$ws.Range("E95").Value = 832
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 588

# Row 112 - Sri Lanka
$ws.Range("D112").Value = 1863
$ws.Range("E112").Value = 192

# Row 114 - Estonia
$ws.Range("B114").Value = 1991
$ws.Range("C114").Value = 1
$ws.Range("D114").Value = 1859
$ws.Range("E114").Value = 63

# Row 132 - Letonia
$ws.Range("D132").Value = 997
$ws.Range("E132").Value = 95

# Row 138 - Georgia
$ws.Range("B138").Value = 943
$ws.Range("C138").Value = 4
$ws.Range("D138").Value = 821

# Row 157 - Taiwan
$ws.Range("B157").Value = 449
$ws.Range("C157").Value = 1
$ws.Range("E157").Value = 4
